$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 10 extra sample rows (rows 14-23) so only 12 data rows remain (2-13)
$ws.Rows("14:23").Delete()

# Update sample names (column A) to the new batch
$ws.Range("A2").Value = "20210426_Cre150_1_20"
$ws.Range("A3").Value = "20210426_Cre150_2_21"
$ws.Range("A4").Value = "20210426_Cre150_3_22"
$ws.Range("A5").Value = "20210426_Cre750_1_28"
$ws.Range("A6").Value = "20210426_Cre750_2_29"
$ws.Range("A7").Value = "20210426_Cre750_3_30"
$ws.Range("A8").Value = "20210426_WT150_1_16"
$ws.Range("A9").Value = "20210426_WT150_2_17"
$ws.Range("A10").Value = "20210426_WT150_3_18"
$ws.Range("A11").Value = "20210426_WT750_1_24"
$ws.Range("A12").Value = "20210426_WT750_2_25"
$ws.Range("A13").Value = "20210426_WT750_3_26"

# Update Resuspension_Volume (B) and Norm1 (D) values, and Norm1_Unit (E) label
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 2).Value = 100
    $ws.Cells.Item($r, 4).Value = 100000
    $ws.Cells.Item($r, 5).Value = "cell"
}

# Recompute bestFit column widths now that content changed
$ws.Columns.Item("A:E").AutoFit()
$ws.Columns.Item(4).ColumnWidth = 7

# Restore plain single-cell selection like the saved file
$ws.Range("H8").Select()
